$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Refresh the panel-query timestamps (time_taken column) on the "data" sheet
$ws.Range("F2").Value = "2021-10-05 14:20:53.860217"
$ws.Range("F3").Value = "2021-10-05 14:20:53.860225"
$ws.Range("F4").Value = "2021-10-05 14:20:53.860228"

# Add a new "metadata" sheet right after "data" to carry the panel query info
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Match the bold/bordered header style used on the "data" sheet
$ws.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

# Data row
$meta.Range("A2").Value = 0
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B2").Value = "Hypocalciuric hypercalcaemia"
$meta.Range("C2").Value = 481
# data_version must stay text ("2.9"), not be coerced to a number
$meta.Range("D2").Value = "'2.9"
$meta.Range("D2").Style = "Normal"
$meta.Range("E2").Value = "2021-07-06T10:53:39.963833Z"
$meta.Range("F2").Value = "2021-10-05 14:20:53.856447"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/481/?format=json"

# Keep "data" as the active sheet/tab, as in the original workbook
$ws.Activate()
